# Trade #47 closed at 2026-02-16 21:30:12 - leadlag DOWN +0.000%
#
# This script updates the "live_trading_results" workbook to reflect:
#   - Trade #18 (leadlag sheet row 17 / All Trades sheet row 18) closing out
#   - A brand-new OPEN trade #47 appended to the leadlag sheet
#   - Updated aggregate stats on the Summary and Comparison sheets

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without Excel's automatic
# type-sniffing turning date-, time- or number-looking strings into real
# dates/numbers. We force the cell to Text format, assign the string, then
# strip the explicit formatting back off again so no stray style survives
# (matches the source file, which carries no cell-level styles at all).
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Summary sheet - OVERALL and leadlag aggregate rows
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 18
Set-TextValue $wsSummary.Range("D2") "61.1%"
Set-TextValue $wsSummary.Range("E2") "+2.1506%"
Set-TextValue $wsSummary.Range("F2") "+0.1195%"

$wsSummary.Range("C3").Value = 35
Set-TextValue $wsSummary.Range("D3") "28.6%"
Set-TextValue $wsSummary.Range("E3") "+2.1053%"
Set-TextValue $wsSummary.Range("F3") "+0.0602%"

# ---------------------------------------------------------------------
# 2) leadlag sheet - close out trade #18 (row 17) and append trade #47
# ---------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

$wsLeadlag.Range("G17").Value = 69028.403791
Set-TextValue $wsLeadlag.Range("H17") "CLOSED"
$wsLeadlag.Range("I17").Value = 0.2673
$wsLeadlag.Range("J17").Value = 2.67
Set-TextValue $wsLeadlag.Range("M17") "time_exit_5min"
$wsLeadlag.Range("N17").Value = 5

$wsLeadlag.Range("A37").Value = 47
Set-TextValue $wsLeadlag.Range("B37") "2026-02-16"
Set-TextValue $wsLeadlag.Range("C37") "21:30:12"
Set-TextValue $wsLeadlag.Range("D37") "leadlag"
Set-TextValue $wsLeadlag.Range("E37") "DOWN"
$wsLeadlag.Range("F37").Value = 68639.325
Set-TextValue $wsLeadlag.Range("H37") "OPEN"
$wsLeadlag.Range("I37").Value = 0
$wsLeadlag.Range("J37").Value = 0
$wsLeadlag.Range("K37").Value = 0.75
Set-TextValue $wsLeadlag.Range("L37") "Coinbase leading with -0.108% move"
$wsLeadlag.Range("N37").Value = 0

# ---------------------------------------------------------------------
# 3) All Trades sheet - append the same closed trade #18 as a new row
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Range("A19").Value = 18
Set-TextValue $wsAll.Range("B19") "2026-02-16"
Set-TextValue $wsAll.Range("C19") "21:25:06"
Set-TextValue $wsAll.Range("D19") "leadlag"
Set-TextValue $wsAll.Range("E19") "DOWN"
$wsAll.Range("F19").Value = 69213.42
$wsAll.Range("G19").Value = 69028.403791
Set-TextValue $wsAll.Range("H19") "CLOSED"
$wsAll.Range("I19").Value = 0.2673
$wsAll.Range("J19").Value = 2.67
$wsAll.Range("K19").Value = 0.75
Set-TextValue $wsAll.Range("L19") "Binance leading with -0.103% move"
Set-TextValue $wsAll.Range("M19") "time_exit_5min"
$wsAll.Range("N19").Value = 5

# ---------------------------------------------------------------------
# 4) Comparison sheet - leadlag strategy row
# ---------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

$wsComparison.Range("B2").Value = 35
Set-TextValue $wsComparison.Range("C2") "28.6%"
Set-TextValue $wsComparison.Range("D2") "1.91"
Set-TextValue $wsComparison.Range("E2") "+0.4418%"
Set-TextValue $wsComparison.Range("G2") "1.15"
